$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells keep their original text representation
# (values like "1.00" / "0.999" must remain text, not be coerced to numbers)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '55.606.25'
$ws.Range("E2").Value = '  -3.25%  '

$ws.Range("D3").Value = '2.906.89'
$ws.Range("E3").Value = '  -3.90%  '

$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '498.50'
$ws.Range("E5").Value = '  -2.95%  '

$ws.Range("D6").Value = '131.96'
$ws.Range("E6").Value = '  -5.44%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("E8").Value = '  -4.22%  '

$ws.Range("D9").Value = '7.11'
$ws.Range("E9").Value = '  -6.04%  '

$ws.Range("E10").Value = '  -6.17%  '

$ws.Range("D11").Value = '0.349'
$ws.Range("E11").Value = '  -5.19%  '

$ws.Range("D12").Value = '3.410.19'
$ws.Range("E12").Value = '  -3.81%  '

$ws.Range("E13").Value = '  -4.25%  '

$ws.Range("D14").Value = '25.54'
$ws.Range("E14").Value = '  -4.22%  '

$ws.Range("D15").Value = '0.0000158'
$ws.Range("E15").Value = '  -4.50%  '

$ws.Range("D16").Value = '55.490.18'
$ws.Range("E16").Value = '  -3.59%  '

$ws.Range("D17").Value = '5.93'
$ws.Range("E17").Value = '  -4.97%  '

$ws.Range("D18").Value = '2.907.53'
$ws.Range("E18").Value = '  -4.08%  '

$ws.Range("D19").Value = '12.57'
$ws.Range("E19").Value = '  -2.08%  '

$ws.Range("D20").Value = '7.63'
$ws.Range("E20").Value = '  -4.85%  '

$ws.Range("D21").Value = '311.87'
$ws.Range("E21").Value = '  -5.97%  '

$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("E23").Value = '  -3.22%  '

$ws.Range("D24").Value = '62.71'
$ws.Range("E24").Value = '  -3.13%  '

$ws.Range("D25").Value = '3.030.76'
$ws.Range("E25").Value = '  -3.87%  '

$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.27%  '

$ws.Range("E27").Value = '  -5.66%  '

$ws.Range("D28").Value = '0.0₃0824'
$ws.Range("E28").Value = '  -11.31%  '

$ws.Range("D29").Value = '6.25'
$ws.Range("E29").Value = '  -8.10%  '

$ws.Range("D30").Value = '6.73'
$ws.Range("E30").Value = '  -10.27%  '

$ws.Range("E31").Value = '  -3.91%  '

$ws.Range("D32").Value = '19.72'
$ws.Range("E32").Value = '  -4.80%  '

$ws.Range("E33").Value = '  -6.69%  '

$ws.Range("D34").Value = '151.74'
$ws.Range("E34").Value = '  -2.33%  '

$ws.Range("D35").Value = '4.34'
$ws.Range("E35").Value = '  -8.47%  '

$ws.Range("D36").Value = '5.56'
$ws.Range("E36").Value = '  -5.65%  '

$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D37").Value = '23.56'
$ws.Range("E37").Value = '  -4.37%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '1.17'
$ws.Range("E38").Value = '  -8.51%  '

$ws.Range("D39").Value = '0.0639'
$ws.Range("E39").Value = '  -6.63%  '

$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '36.38'
$ws.Range("E40").Value = '  -2.89%  '

$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.12%  '

$ws.Range("D42").Value = '3.67'
$ws.Range("E42").Value = '  -5.11%  '

$ws.Range("E43").Value = '  -3.01%  '

$ws.Range("D44").Value = '5.96'
$ws.Range("E44").Value = '  -1.60%  '

$ws.Range("D45").Value = '2.098.43'
$ws.Range("E45").Value = '  -8.93%  '

$ws.Range("E46").Value = '  -6.61%  '

$ws.Range("D47").Value = '0.910'
$ws.Range("E47").Value = '  -8.14%  '

$ws.Range("D48").Value = '0.0232'
$ws.Range("E48").Value = '  -3.25%  '

$ws.Range("D49").Value = '18.46'
$ws.Range("E49").Value = '  -5.79%  '

$ws.Range("D50").Value = '0.0836'
$ws.Range("E50").Value = '  -6.42%  '

$ws.Range("D51").Value = '1.64'
$ws.Range("E51").Value = '  -11.39%  '
